# Round specific statistic values to 4 decimal places on the
# "normality" and "pairwise_ttests" sheets.

$wb = $excel.ActiveWorkbook

$wsNormality = $wb.Worksheets.Item("normality")
$wsNormality.Range("B3").Value = 0.9654
$wsNormality.Range("C3").Value = 0.8575
$wsNormality.Range("B4").Value = 0.8455
$wsNormality.Range("C4").Value = 0.025

$wsPairwise = $wb.Worksheets.Item("pairwise_ttests")
$wsPairwise.Range("G3").Value = -0.0668
$wsPairwise.Range("H3").Value = 22.8676
$wsPairwise.Range("J3").Value = 0.9473
$wsPairwise.Range("L3").Value = -0.0257
